$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for each data row (2..82).
# This update bumps that date forward by one day (46060 -> 46061) for every row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 82 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cur = $cell.Value2
    if ($cur -ne $null) {
        $cell.Value2 = $cur + 1
    }
}
